$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question text in C2 (was "Haben Sie ein Auto?", now "Haben Sie ein Auto??")
$ws.Range("C2").Value = "Haben Sie ein Auto??"

# Update selection / view: active cell moves from E3 to C2, and the sheet is
# scrolled back so column A is the left-most visible column again.
$ws.Range("A1").Select() | Out-Null
$ws.Range("C2").Select() | Out-Null
